$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.063.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.58%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.70%  "

# Row 4
$ws.Range("E4").Value = "  +0.42%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.618"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.62%  "

# Row 7
$ws.Range("E7").Value = "  +0.36%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.78%  "

# Row 9
$ws.Range("E9").Value = "  +3.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0694"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.51%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0984"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.74%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.121.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.72%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.68%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.676"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.03%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.839.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.94%  "

# Row 16
$ws.Range("E16").Value = "  +3.29%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.104.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.73%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.76%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0793"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.02%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.53%  "

# Row 23
$ws.Range("E23").Value = "  +0.37%  "

# Row 24
$ws.Range("E24").Value = "  +1.98%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "173.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "

# Row 26
$ws.Range("E26").Value = "  +1.91%  "

# Row 27
$ws.Range("E27").Value = "  +2.18%  "

# Row 28
$ws.Range("E28").Value = "  +4.41%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.73%  "

# Row 30
$ws.Range("E30").Value = "  +0.34%  "

# Row 31
$ws.Range("E31").Value = "  +1.73%  "

# Row 32
$ws.Range("E32").Value = "  +0.07%  "

# Row 33
$ws.Range("E33").Value = "  +1.96%  "

# Row 34
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.79%  "

# Row 35
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +22.61%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.764"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.60%  "

# Row 37
$ws.Range("E37").Value = "  +7.74%  "

# Row 38
$ws.Range("E38").Value = "  +13.43%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "90.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.45%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.353.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.81%  "

# Row 41
$ws.Range("E41").Value = "  +3.21%  "

# Row 42
$ws.Range("E42").Value = "  +1.90%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.26%  "

# Row 44
$ws.Range("E44").Value = "  -1.53%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.45%  "

# Row 46
$ws.Range("E46").Value = "  +4.47%  "

# Row 47
$ws.Range("E47").Value = "  +3.92%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.040.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.79%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +19.24%  "

# Row 50
$ws.Range("E50").Value = "  +0.39%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0669"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
